# Update the dSF column (F) values to reflect repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    3  = 0
    4  = -3
    6  = -2
    9  = -2
    11 = -3
    13 = 1
    20 = -3
    21 = -3
    22 = -2
    28 = -4
    31 = 7
    34 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
